$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Grading corrections (column D) for several students
# ---------------------------------------------------------------
$ws.Range("D2").Value  = 5.75   # Sharma Kartik
$ws.Range("D7").Value  = 5.75   # Javaid Mohammad
$ws.Range("D8").Value  = 5      # Sen Majumder Riddhi
$ws.Range("D9").Value  = 5.75   # Sokolowska Julia Magdalena
$ws.Range("D10").Value = 5.75   # Biben Valeria
$ws.Range("D15").Value = 5.75   # Garonzi Margherita
$ws.Range("D17").Value = 5.75   # Ramaswamy Krithi Dakshina

# ---------------------------------------------------------------
# 2. Insert a blank separator row between the regular students
#    (rows 2-19) and the auditing students (rows 20 onward).
# ---------------------------------------------------------------
$ws.Rows.Item(20).Insert()

# style/format the new blank row 20 like the other data rows, and
# match its (slightly shorter) row height
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = $null
$ws.Rows.Item(20).RowHeight = 17

# ---------------------------------------------------------------
# 3. "Dhananjayan Parvati" (row 19) is confirmed to be auditing, so
#    she moves down below the separator (row 21, which currently
#    holds "Manish Medha", shifted down by the insert above).
#    "Manish Medha" is confirmed to be attending normally (not
#    auditing) and takes row 19 instead. Swap the two names via a
#    scratch cell (C2 is blank and already inside the used range)
#    so the existing shared-string entries are reused as-is.
# ---------------------------------------------------------------
$ws.Range("A19").Copy()
$ws.Range("C2").PasteSpecial(-4104)
$ws.Range("A21").Copy()
$ws.Range("A19").PasteSpecial(-4104)
$ws.Range("C2").Copy()
$ws.Range("A21").PasteSpecial(-4104)
$ws.Range("C2").ClearContents()

# Row 19 -> "Manish Medha", attending normally, full grade.
$ws.Range("B19").Value = 0
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 1
$ws.Range("F19").ClearContents()

# Row 21 -> "Dhananjayan Parvati", confirmed auditing.
$ws.Range("B21").Value = 1
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").ClearContents()

# ---------------------------------------------------------------
# 4. Restore the active selection as left by the editor
# ---------------------------------------------------------------
$ws.Range("E12").Select()
